$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark task 29 (row 32) as Done, add time estimate
$ws.Range("F32").Value = "Done"
$ws.Range("F32").Interior.Color = 5296274
$ws.Range("F32").WrapText = $true
$ws.Range("G32").Value = "1 hour"

# Mark task 30 (row 33) as Done, add time estimate
$ws.Range("F33").Value = "Done"
$ws.Range("F33").Interior.Color = 5296274
$ws.Range("F33").WrapText = $true
$ws.Range("G33").Value = "1 hour"

# Update selection to reflect the last edited cell (F33)
$ws.Range("F33").Select()
